$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 135, shifting the existing rows 135-227 down to 136-228.
$ws.Rows.Item(135).Insert()

# Populate the newly inserted row 135 with the new weekly record.
$ws.Cells.Item(135, 1).Value = 7
$ws.Cells.Item(135, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(135, 3).Value = "Ñuble"
$ws.Cells.Item(135, 4).Value = 44603
$ws.Cells.Item(135, 5).Value = 16
$ws.Cells.Item(135, 6).Value = 100112008
$ws.Cells.Item(135, 7).Value = "Coliflor"
$ws.Cells.Item(135, 8).Value = "Sin especificar"
$ws.Cells.Item(135, 9).Value = "Primera"
$ws.Cells.Item(135, 10).Value = 200
$ws.Cells.Item(135, 11).Value = 800
$ws.Cells.Item(135, 12).Value = 850
$ws.Cells.Item(135, 13).Value = 825
$ws.Cells.Item(135, 14).Value = "$/unidad"
$ws.Cells.Item(135, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(135, 16).Value = 825
$ws.Cells.Item(135, 17).Value = 1
$ws.Cells.Item(135, 18).Value = "Hortaliza"

# Make sure the new date cell keeps the same date number-format style used by
# the rest of column D.
$ws.Cells.Item(135, 4).NumberFormat = $ws.Cells.Item(136, 4).NumberFormat
